$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing rows store every non-jobId column as literal text (the
# jobId column A is the only numeric column). Force each text cell to
# "@" (Text) format right before assigning so values like "90", "900" or
# "2019-12-27" are kept as strings instead of being auto-converted to
# numbers or dates.

# Row 10 (jobId 9)
$ws.Range("A10").Value = 9
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "david"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "2019-12-27"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "pmma"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "Cut"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "90"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "90"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "900"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "5000"
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = "1/0"
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = "1"
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "Enter here useful comments for the future"

# Row 11 (jobId 10)
$ws.Range("A11").Value = 10
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "david"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "2019-12-27"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "pmma"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "Cut"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "90"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "90"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "900"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "5000"
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = "1/0"
$ws.Range("K11").NumberFormat = "@"
$ws.Range("K11").Value = "1"
$ws.Range("L11").NumberFormat = "@"
$ws.Range("L11").Value = "Enter here useful comments for the future"

# Row 12 (jobId 11)
$ws.Range("A12").Value = 11
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "david"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "2019-12-27"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "pmma"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "Cut"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "90"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "90"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "900"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "5000"
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = "1/0"
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value = "1"
$ws.Range("L12").NumberFormat = "@"
$ws.Range("L12").Value = "Enter here useful comments for the future"

# Row 13 (jobId 12)
$ws.Range("A13").Value = 12
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "david"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "2019-12-27"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "pmma"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "Cut"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "90"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "90"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "900"
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "5000"
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value = "1/0"
$ws.Range("K13").NumberFormat = "@"
$ws.Range("K13").Value = "1"
$ws.Range("L13").NumberFormat = "@"
$ws.Range("L13").Value = "Enter here useful comments for the future"

# Row 14 (jobId 13)
$ws.Range("A14").Value = 13
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "2019-12-27"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "Cut"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "90"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "90"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "900"
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "5000"
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = "1/0"
$ws.Range("K14").NumberFormat = "@"
$ws.Range("K14").Value = "1"
$ws.Range("L14").NumberFormat = "@"
$ws.Range("L14").Value = "Enter here useful comments for the future"

# Row 15 (jobId 14)
$ws.Range("A15").Value = 14
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "David"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "2019-12-27"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "pmma"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "Cut"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "90"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "90"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "900"
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "5000"
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = "1/0"
$ws.Range("K15").NumberFormat = "@"
$ws.Range("K15").Value = "1"
$ws.Range("L15").NumberFormat = "@"
$ws.Range("L15").Value = "Enter here useful comments for the future"

# Row 16 (jobId 15)
$ws.Range("A16").Value = 15
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "david"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "2019-12-27"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "sache"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "Cut"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "90"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "90"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "900"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "5000"
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "1/0"
$ws.Range("K16").NumberFormat = "@"
$ws.Range("K16").Value = "1"
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = "Enter here useful comments for the future"

# Row 17 (jobId 16)
$ws.Range("A17").Value = 16
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "hola"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "2019-12-27"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "hola"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "Cut"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "90"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "90"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "900"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "5000"
$ws.Range("J17").NumberFormat = "@"
$ws.Range("J17").Value = "1/0"
$ws.Range("K17").NumberFormat = "@"
$ws.Range("K17").Value = "1"
$ws.Range("L17").NumberFormat = "@"
$ws.Range("L17").Value = "Enter here useful comments for the future"

# Row 18 (jobId 17)
$ws.Range("A18").Value = 17
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "hola"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "2019-12-27"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "hola"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "Cut"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "90"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "90"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "900"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "5000"
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = "1/0"
$ws.Range("K18").NumberFormat = "@"
$ws.Range("K18").Value = "1"
$ws.Range("L18").NumberFormat = "@"
$ws.Range("L18").Value = "Enter here useful comments for the future"

# Row 19 (jobId 18)
$ws.Range("A19").Value = 18
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "2019-12-27"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "Cut"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "90"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "90"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "900"
$ws.Range("I19").NumberFormat = "@"
$ws.Range("I19").Value = "5000"
$ws.Range("J19").NumberFormat = "@"
$ws.Range("J19").Value = "1/0"
$ws.Range("K19").NumberFormat = "@"
$ws.Range("K19").Value = "1"
$ws.Range("L19").NumberFormat = "@"
$ws.Range("L19").Value = "Enter here useful comments for the future"

# Row 20 (jobId 19)
$ws.Range("A20").Value = 19
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = ""
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "2019-12-27"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = ""
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "Cut"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "90"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "90"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "900"
$ws.Range("I20").NumberFormat = "@"
$ws.Range("I20").Value = "5000"
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "1/0"
$ws.Range("K20").NumberFormat = "@"
$ws.Range("K20").Value = "1"
$ws.Range("L20").NumberFormat = "@"
$ws.Range("L20").Value = "Enter here useful comments for the future"
